$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The topic in cell D11 ("Data pipelines") is being renamed to
# "Data pipelines and reproducibility".
$ws.Range("D11").Value = "Data pipelines and reproducibility"

# Reflect the active cell selection moving to D12 (matches the diff's
# sheetView selection change).
$ws.Range("D12").Select()
